$d = $word.ActiveDocument

# Locate the placeholder run that needs to be replaced/split.
$r = $d.Content
$r.Find.Execute("#<CompanyDetails_Type>#")
$start = $r.Start
$end = $r.End

# Replace the whole placeholder text with the new combined text
# (this keeps the original run formatting - Times New Roman, sz 24).
$newText = "Oświadczam, że prowadzę #<InvestorType_B>#"
$r.Text = $newText

# Compute the boundaries of the four resulting runs:
#   1) "Oświadczam, że prowadzę "   (24 chars)
#   2) "#<"                        (2 chars)
#   3) "InvestorType_B"            (14 chars)
#   4) ">#"                        (2 chars)
$b0 = $start
$b1 = $start + 24
$b2 = $start + 26
$b3 = $start + 40
$b4 = $start + 42

# Force Word to split the single run into separate runs at each boundary
# by toggling a character formatting property on/off over each segment
# (this leaves the visible formatting unchanged but creates run breaks).
$seg2 = $d.Range($b1, $b2)
$seg2.Bold = 1
$seg2.Bold = 0

$seg3 = $d.Range($b2, $b3)
$seg3.Bold = 1
$seg3.Bold = 0

$seg4 = $d.Range($b3, $b4)
$seg4.Bold = 1
$seg4.Bold = 0
